$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition) - column F "想去人数" (want-to-go count) updates
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 32
$ws1.Range("F5").Value = 5250
$ws1.Range("F6").Value = 5250
$ws1.Range("F7").Value = 152
$ws1.Range("F9").Value = 532
$ws1.Range("F13").Value = 5184
$ws1.Range("F15").Value = 73
$ws1.Range("F16").Value = 91
$ws1.Range("F17").Value = 290
$ws1.Range("F18").Value = 290
$ws1.Range("F19").Value = 255
$ws1.Range("F22").Value = 3904
$ws1.Range("F24").Value = 3827
$ws1.Range("F29").Value = 247
$ws1.Range("F30").Value = 208
$ws1.Range("F36").Value = 19
$ws1.Range("F37").Value = 6803
$ws1.Range("F38").Value = 1103
$ws1.Range("F42").Value = 60
$ws1.Range("F43").Value = 1385
$ws1.Range("F44").Value = 173
$ws1.Range("F45").Value = 692
$ws1.Range("F48").Value = 314
$ws1.Range("F49").Value = 92

# Sheet 2: 演出 (Performance) - column F updates
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F11").Value = 65
$ws2.Range("F17").Value = 142

# Sheet 2 row 5: column G "最低票价" (lowest ticket price) -> now unavailable for sale
$ws2.Range("G5").Value = "不可售"

# Sheet 4: 全部类型 (All types) - column F updates
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F7").Value = 5250
$ws4.Range("F8").Value = 5250
$ws4.Range("F9").Value = 152
$ws4.Range("F12").Value = 532
$ws4.Range("F15").Value = 5184
$ws4.Range("F17").Value = 73
$ws4.Range("F18").Value = 91
$ws4.Range("F19").Value = 290
$ws4.Range("F20").Value = 290
$ws4.Range("F21").Value = 255
$ws4.Range("F24").Value = 3904
$ws4.Range("F25").Value = 3827
$ws4.Range("F29").Value = 247
$ws4.Range("F30").Value = 208
$ws4.Range("F35").Value = 19
$ws4.Range("F36").Value = 142
$ws4.Range("F37").Value = 6803
$ws4.Range("F38").Value = 1103
$ws4.Range("F43").Value = 60
$ws4.Range("F44").Value = 1385
$ws4.Range("F45").Value = 173
$ws4.Range("F46").Value = 692
$ws4.Range("F48").Value = 314
